# Applies the "Fix bug in prepare pipe / Refactor to support nlp load from
# base / Update data set / Update README" edit to the workbook:
#  - rename "train" sheet to "source" and collapse it to a single summary
#    row (idx/kind/text)
#  - insert a brand new "prepare" sheet (phrase-level NER data) between
#    "source" and "config"
#  - update the "config" sheet keys/values (dotted names -> underscore
#    names, add nlp_base/source_language/prepare_enabled/train_autosave)
#  - make "config" the active/selected sheet
#
# NOTE: worksheet handles in this COM runtime are positional - inserting
# sheets (Worksheets.Add) re-indexes previously obtained handles for any
# sheet whose position shifts, so we re-fetch a worksheet by name
# immediately before using it whenever a structural operation could have
# moved it. Also, the shared-string table is built in first-write order,
# so cell writes below are intentionally sequenced to reproduce the
# target's string order.

$wb = $excel.ActiveWorkbook

$xlLeft = -4131

# ---------------------------------------------------------------------
# 1. "train" -> "source" (still at position 1 throughout this script)
# ---------------------------------------------------------------------
$wsSource = $wb.Worksheets.Item("train")
$wsSource.Name = "source"
$wsSource.Cells.Clear()

$wsSource.Range("A1").Value = "idx"
$wsSource.Range("B1").Value = "kind"
$wsSource.Range("C1").Value = "text"

$wsSource.Columns.Item(1).ColumnWidth = 6.5
$wsSource.Columns.Item(2).ColumnWidth = 8.33203125
$wsSource.Columns.Item(3).ColumnWidth = 52.83203125
$wsSource.Columns.Item(4).ColumnWidth = 15.1640625

# ---------------------------------------------------------------------
# 2. brand new "prepare" sheet, inserted right before "config" (i.e.
#    right after "source"); re-fetch "config" immediately before the
#    insertion since it is the sheet whose position is about to shift.
# ---------------------------------------------------------------------
$wsConfigAnchor = $wb.Worksheets.Item("config")
$wsPrepare = $wb.Worksheets.Add($wsConfigAnchor)
$wsPrepare.Name = "prepare"

$wsPrepare.Range("A1").Value = "idx"
$wsPrepare.Range("B1").Value = "kind"
$wsPrepare.Range("C1").Value = "value"
$wsPrepare.Range("D1").Value = "entity"

$wsPrepare.Range("A2").Value = 1
$wsPrepare.Range("B2").Value = "phrase"
$wsPrepare.Range("C2").Value = "Robertus Johansyah"
$wsPrepare.Range("D2").Value = "PERSON"

$wsPrepare.Range("A3").Value = 2
$wsPrepare.Range("B3").Value = "phrase"
$wsPrepare.Range("C3").Value = "ExcelCy"
$wsPrepare.Range("D3").Value = "GITHUB_PROJECT"

$wsPrepare.Range("A4").Value = 3
$wsPrepare.Range("B4").Value = "phrase"
$wsPrepare.Range("C4").Value = "Github"
$wsPrepare.Range("D4").Value = "PRODUCT"

$wsPrepare.Columns.Item(3).ColumnWidth = 19.1640625
$wsPrepare.Columns.Item(4).ColumnWidth = 17.1640625

$wsPrepare.Range("C3").Select()

# ---------------------------------------------------------------------
# 3. back to "source" (handle is still valid - its position never
#    moved) to write the data row, after "prepare"'s strings so the
#    shared-string table order matches the target.
# ---------------------------------------------------------------------
$wsSource.Range("A2").Value = 1
$wsSource.Range("B2").Value = "text"
$wsSource.Range("C2").Value = "Robertus Johansyah is the maintainer of project ExcelCy in Github"

# ---------------------------------------------------------------------
# 4. "config" sheet: new keys/values, new row, becomes the active tab.
#    Re-fetch by name since its position shifted when "prepare" was
#    inserted in front of it.
# ---------------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("config")

$wsConfig.Range("A1").Value = "name"
$wsConfig.Range("B1").Value = "value"

$wsConfig.Range("A2").Value = "nlp_base"
$wsConfig.Range("B2").Value = "en_core_web_sm"

$wsConfig.Range("A3").Value = "nlp_name"
$wsConfig.Range("B3").Value = "[tmp]/nlp/test_data_28"

$wsConfig.Range("A4").Value = "source_language"
$wsConfig.Range("B4").Value = "en"

$wsConfig.Range("A5").Value = "prepare_enabled"
$wsConfig.Range("B5").Value = $true
$wsConfig.Range("B5").HorizontalAlignment = $xlLeft

$wsConfig.Range("A6").Value = "train_iteration"
$wsConfig.Range("B6").Value = 2

$wsConfig.Range("A7").Value = "train_drop"
$wsConfig.Range("B7").Value = 0.2

$wsConfig.Range("A8").Value = "train_autosave"
$wsConfig.Range("B8").Value = $true
$wsConfig.Range("B8").HorizontalAlignment = $xlLeft

$wsConfig.Activate()
